$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated player pool data (players reshuffled / refreshed for row 2-19)
$data = @(
    @("Derrick White",        "PG,SG",    "Boston Celtics"),
    @("Collin Sexton",        "PG,SG",    "Utah Jazz"),
    @("Jrue Holiday",         "PG,SG",    "Boston Celtics"),
    @("Max Christie",         "SG,SF",    "Los Angeles Lakers"),
    @("Herbert Jones",        "SF,PF",    "New Orleans Pelicans"),
    @("Anthony Davis",        "PF,C",     "Los Angeles Lakers"),
    @("Malik Monk",           "PG,SG,SF", "Sacramento Kings"),
    @("Yves Missi",           "C",        "New Orleans Pelicans"),
    @("Isaiah Hartenstein",   "C",        "Oklahoma City Thunder"),
    @("Bam Adebayo",          "C",        "Miami Heat"),
    @("Scotty Pippen Jr.",    "PG,SG",    "Memphis Grizzlies"),
    @("Damian Lillard",       "PG",       "Milwaukee Bucks"),
    @("Cade Cunningham",      "PG,SG",    "Detroit Pistons"),
    @("Donte DiVincenzo",     "SG,SF",    "Minnesota Timberwolves"),
    @("Julius Randle",        "PF,C",     "Minnesota Timberwolves"),
    @("Brandon Miller",       "SG,SF,PF", "Charlotte Hornets"),
    @("LaMelo Ball",          "PG,SG",    "Charlotte Hornets"),
    @("Cameron Johnson",      "SF,PF",    "Brooklyn Nets")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}
